$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the trailing "Resolving-Mac" sending-cluster rows (old rows 14:17) --
# the target cluster "Resolving-Mac" was dropped from the TPM re-run, so the
# sending-cluster's own "Resolving-Mac" block (senderless of the Resolving-Mac
# target) no longer applies to this ligand-receptor pair's recomputed table.
$ws.Rows("14:17").Delete()

# Recomputed TPM-based values (updated ligand/receptor expression & specificity
# figures, and target cluster "Resolving-Mac" removed as a receiving cluster)
# for the remaining sending/target cluster combinations.
$data = @(
  @("ECs","Wnt2b","Fzd4","ECs",3,1,1.993530333333333,5.980591,0.2657535748877558,0.2657535748877558,3,1,30.58864766666666,91.76594299999999,0.3925391465174898,0.3925391465174898,60.97939697914587,548.8145728123129,0.1043186814704115,0.1043186814704114),
  @("ECs","Wnt2b","Fzd4","FAPs",3,1,1.993530333333333,5.980591,0.2657535748877558,0.2657535748877558,3,1,22.690535,68.07160500000001,0.291183949679193,0.291183949679193,45.23426980206166,407.108428218555,0.07738317557718194,0.07738317557718191),
  @("ECs","Wnt2b","Fzd4","MuSCs",3,1,1.993530333333333,5.980591,0.2657535748877558,0.2657535748877558,3,1,24.64590566666666,73.93771699999999,0.3162769038033173,0.3162769038033172,49.13236053897188,442.1912448507469,0.08405171784016242,0.0840517178401624),
  @("FAPs","Wnt2b","Fzd4","ECs",3,1,2.433537333333333,7.300612,0.324410035374171,0.324410035374171,3,1,30.58864766666666,91.76594299999999,0.3925391465174898,0.3925391465174898,74.43861607301288,669.9475446571159,0.1273436384074857,0.1273436384074857),
  @("FAPs","Wnt2b","Fzd4","FAPs",3,1,2.433537333333333,7.300612,0.324410035374171,0.324410035374171,3,1,22.690535,68.07160500000001,0.291183949679193,0.291183949679193,55.21826403580666,496.96437632226,0.09446299541581785,0.0944629954158178),
  @("FAPs","Wnt2b","Fzd4","MuSCs",3,1,2.433537333333333,7.300612,0.324410035374171,0.324410035374171,3,1,24.64590566666666,73.93771699999999,0.3162769038033173,0.3162769038033172,59.97673155364488,539.7905839828039,0.1026034015508674,0.1026034015508674),
  @("MuSCs","Wnt2b","Fzd4","ECs",3,1,2.135898,6.407693999999999,0.2847323261675683,0.2847323261675683,3,1,30.58864766666666,91.76594299999999,0.3925391465174898,0.3925391465174898,65.33423137393798,588.0080823654419,0.1117685842997568,0.1117685842997568),
  @("MuSCs","Wnt2b","Fzd4","FAPs",3,1,2.135898,6.407693999999999,0.2847323261675683,0.2847323261675683,3,1,22.690535,68.07160500000001,0.291183949679193,0.291183949679193,48.46466832542999,436.18201492887,0.08290948333481678,0.08290948333481676),
  @("MuSCs","Wnt2b","Fzd4","MuSCs",3,1,2.135898,6.407693999999999,0.2847323261675683,0.2847323261675683,3,1,24.64590566666666,73.93771699999999,0.3162769038033173,0.3162769038033172,52.64114062162199,473.7702655945979,0.09005425853299476,0.09005425853299474),
  @("Resolving-Mac","Wnt2b","Fzd4","ECs",3,1,0.9384586666666667,2.815376,0.125104063570505,0.125104063570505,3,1,30.58864766666666,91.76594299999999,0.3925391465174898,0.3925391465174898,28.70618150439644,258.355633539568,0.04910824233983584,0.04910824233983584),
  @("Resolving-Mac","Wnt2b","Fzd4","FAPs",3,1,0.9384586666666667,2.815376,0.125104063570505,0.125104063570505,3,1,22.690535,68.07160500000001,0.291183949679193,0.291183949679193,21.29412922205333,191.64716299848,0.03642829535137651,0.0364282953513765),
  @("Resolving-Mac","Wnt2b","Fzd4","MuSCs",3,1,0.9384586666666667,2.815376,0.125104063570505,0.125104063570505,3,1,24.64590566666666,73.93771699999999,0.3162769038033173,0.3162769038033172,23.12916377073244,208.162473936592,0.03956752587929271,0.03956752587929271)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowArr = $data[$r]
    for ($c = 0; $c -lt $rowArr.Length; $c++) {
        $ws.Cells.Item(2 + $r, 1 + $c).Value = $rowArr[$c]
    }
}
